# Changed results for automotive/* and consumer/jpeg to runme_large.sh command.
#
# The sheet shrinks from 9 data rows (A1:D9) to 6 (A1:D6): basicmath's two
# commands collapse into one row, and bitcount/qsort/susan each collapse
# from their old multi-variant commands (2-3 rows apiece) into a single
# "runme_large.sh" row, plus a new consumer/jpeg row is appended at the end.
#
# Rows 2, 4, 6, 8 (basicmathsmall, bitcnts, qsort_large, susan -e) are the
# ones kept/reused in place (so their original formatting carries through
# correctly); rows 3, 5, 7, 9 are removed. Deleting bottom-up avoids row
# index shifting while we work.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("9:9").Delete()
$ws.Rows("7:7").Delete()
$ws.Rows("5:5").Delete()
$ws.Rows("3:3").Delete()

# Row 2: automotive/basicmath/* -> automotive/basicmath/runme_large.sh
$ws.Range("A2").Value = "automotive/basicmath/runme_large.sh"
$ws.Range("B2").Value = 0.24
$ws.Range("C2").Value = 0.23
$ws.Range("D2").Value = 0.01

# Row 3: automotive/bitcount/* -> automotive/bitcount/runme_large.sh
$ws.Range("A3").Value = "automotive/bitcount/runme_large.sh"
$ws.Range("B3").Value = 0.06
$ws.Range("C3").Value = 0.05
$ws.Range("D3").Value = 0

# Row 4: automotive/qsort/* -> automotive/qsort/runme_large.sh
$ws.Range("A4").Value = "automotive/qsort/runme_large.sh"
$ws.Range("B4").Value = 0.05
$ws.Range("C4").Value = 0.04
$ws.Range("D4").Value = 0.01

# Row 5: automotive/susan/* -> automotive/susan/runme_large.sh
$ws.Range("A5").Value = "automotive/susan/runme_large.sh"
$ws.Range("B5").Value = 0.06
$ws.Range("C5").Value = 0.05
$ws.Range("D5").Value = 0

# Row 6 (new): consumer/jpeg/runme_large.sh
$ws.Range("A6").Value = "consumer/jpeg/runme_large.sh"
$ws.Range("B6").Value = 0.02
$ws.Range("C6").Value = 0.02
$ws.Range("D6").Value = 0

# Match the author's final selection/cursor position.
$ws.Range("A8").Select() | Out-Null
